# Insert a new weekly record for "Agrícola del Norte S.A. de Arica - Zanahoria".
# The new record is inserted as row 138 (pushing the previous rows 138-238
# down to 139-239), keeping Calidad ("Primera") and Origen ("Valle de Camiña")
# the same as the row it is inserted before, but with fresh Fecha / Volumen /
# Precio values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 138..238 down to 139..239, duplicating row 138's formatting.
$ws.Rows.Item(138).Insert()

# Seed the new row 138 with the same static field values as the row that
# used to occupy position 138 (now shifted to 139), since Mercado, Región,
# Categoría, Variedad, Calidad, Unidad de comercialización, Origen, Kg o
# Unidades and Clasificación stay the same for this new record.
$ws.Range("A139:R139").Copy()
$ws.Range("A138").PasteSpecial()

# Now overwrite the fields that actually differ for the new record.
$ws.Range("D138").Value = 44603
$ws.Range("J138").Value = 90
$ws.Range("K138").Value = 17000
$ws.Range("L138").Value = 18000
$ws.Range("M138").Value = 17500
$ws.Range("P138").Value = 700
